$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header shared strings: bump the weekly report volume/number and date range.
# These shared strings are rich text (multiple runs); use Characters() to
# replace just the run's text in place instead of re-writing the whole cell.
# ---------------------------------------------------------------------------

# A8: "Volume 32   Number  49" -> "...50"
$ws.Range("A8").Characters(21, 2).Text = "50"

# C9: "Report Covering the Week  12/1/2025  Through  12/7/2025"
#  -> "...12/8/2025  Through  12/14/2025"
$ws.Range("C9").Characters(27, 9).Text = "12/8/2025"
$ws.Range("C9").Characters(47, 9).Text = "12/14/2025"

# ---------------------------------------------------------------------------
# Helper donor cells (untouched by this edit) used to coerce a cell's style
# when its type flips between a numeric stat and the "no activity" text
# placeholders ("0" / "***.*").
# ---------------------------------------------------------------------------
$textStyleDonor = $ws.Range("C14")   # s=13, general text placeholder style
$numStyleDonor  = $ws.Range("I14")   # s=14, plain number style
$pctStyleDonor  = $ws.Range("K14")   # s=15, percent/decimal number style

function Set-AsText($cell, [string]$text) {
    $textStyleDonor.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
    $cell.Value = $text
}

function Set-AsNumber($cell, $value) {
    $numStyleDonor.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
    $cell.Value = $value
}

function Set-AsPercent($cell, $value) {
    $pctStyleDonor.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
    $cell.Value = $value
}

# Row 14 - Murder
$ws.Range("M14").Value = -44.444444444444
$ws.Range("N14").Value = -86.486486486486

# Row 15 - Rape
$ws.Range("F15").Value = 3
$ws.Range("I15").Value = 19
$ws.Range("K15").Value = 18.75
$ws.Range("L15").Value = 5.555555555555
$ws.Range("M15").Value = 46.153846153846
$ws.Range("N15").Value = -73.239436619718

# Row 16 - Robbery
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 60
$ws.Range("I16").Value = 135
$ws.Range("J16").Value = 85
$ws.Range("K16").Value = 58.823529411764
$ws.Range("L16").Value = 8
$ws.Range("M16").Value = -49.248120300751
$ws.Range("N16").Value = -91.640866873065

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 41
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = 95.238095238095
$ws.Range("I17").Value = 417
$ws.Range("J17").Value = 262
$ws.Range("K17").Value = 59.160305343511
$ws.Range("L17").Value = 20.172910662824
$ws.Range("M17").Value = 45.296167247386
$ws.Range("N17").Value = -53.666666666666

# Row 18 - Burglary
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 160
$ws.Range("I18").Value = 87
$ws.Range("J18").Value = 71
$ws.Range("K18").Value = 22.535211267605
$ws.Range("L18").Value = 12.987012987013
$ws.Range("M18").Value = -70.205479452054
$ws.Range("N18").Value = -95.188053097345

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 200
$ws.Range("F19").Value = 19
$ws.Range("G19").Value = 13
$ws.Range("H19").Value = 46.153846153846
$ws.Range("I19").Value = 330
$ws.Range("J19").Value = 232
$ws.Range("K19").Value = 42.241379310344
$ws.Range("L19").Value = -16.876574307304
$ws.Range("M19").Value = -19.117647058823
$ws.Range("N19").Value = -50.598802395209

# Row 20 - G.L.A.
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("J20").Value = 77
$ws.Range("K20").Value = 27.272727272727
$ws.Range("L20").Value = 6.521739130434
$ws.Range("M20").Value = -32.876712328767
$ws.Range("N20").Value = -92.570128885519

# Row 21 - TOTAL
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = 120
$ws.Range("F21").Value = 99
$ws.Range("G21").Value = 59
$ws.Range("H21").Value = 67.796610169491
$ws.Range("I21").Value = 1091
$ws.Range("J21").Value = 748
$ws.Range("K21").Value = 45.855614973262
$ws.Range("L21").Value = 2.827521206409
$ws.Range("M21").Value = -23.223082336382
$ws.Range("N21").Value = -83.000934870676

# Row 22 - Transit
$ws.Range("M22").Value = -80.434782608695

# Row 23 - Housing (L23 flips from the "***.*" text placeholder to a number)
Set-AsPercent $ws.Range("L23") 100

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 19
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 72
$ws.Range("G24").Value = 80
$ws.Range("H24").Value = -10
$ws.Range("I24").Value = 915
$ws.Range("J24").Value = 864
$ws.Range("K24").Value = 5.902777777777
$ws.Range("L24").Value = -14.804469273743
$ws.Range("M24").Value = -5.864197530864

# Row 25 - Retail Theft (D25, E25 flip to text placeholders)
Set-AsText $ws.Range("D25") "0"
Set-AsText $ws.Range("E25") "***.*"
$ws.Range("F25").Value = 9
$ws.Range("H25").Value = -18.181818181818
$ws.Range("I25").Value = 190
$ws.Range("K25").Value = 29.251700680272
$ws.Range("L25").Value = -19.148936170212

# Row 26 - Misd. Assault
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = -28.571428571428
$ws.Range("F26").Value = 44
$ws.Range("G26").Value = 46
$ws.Range("H26").Value = -4.347826086956
$ws.Range("I26").Value = 488
$ws.Range("J26").Value = 511
$ws.Range("K26").Value = -4.500978473581
$ws.Range("L26").Value = -0.204498977505
$ws.Range("M26").Value = -38.770388958594

# Row 27 - UCR Rape*
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 19
$ws.Range("K27").Value = -20.833333333333
$ws.Range("L27").Value = -13.636363636363

# Row 28 - Other Sex Crimes (C28 flips from the "0" text placeholder to a number)
Set-AsNumber $ws.Range("C28") 1
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -33.333333333333
$ws.Range("I28").Value = 55
$ws.Range("J28").Value = 64
$ws.Range("K28").Value = -14.0625
$ws.Range("L28").Value = -1.785714285714

# Row 29 - Shooting Vic. (C29, D29, E29 flip to text placeholders)
Set-AsText $ws.Range("C29") "0"
Set-AsText $ws.Range("D29") "0"
Set-AsText $ws.Range("E29") "***.*"
$ws.Range("N29").Value = -71.578947368421

# Row 30 - Shooting Inc. (C30, D30, E30 flip to text placeholders)
Set-AsText $ws.Range("C30") "0"
Set-AsText $ws.Range("D30") "0"
Set-AsText $ws.Range("E30") "***.*"
$ws.Range("N30").Value = -88.172043010752

$excel.CutCopyMode = 0
